$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 24,14
$data[0,0] = 8.063151185903664
$data[0,1] = 4.584180610710706
$data[0,2] = 0
$data[0,3] = 12.74416604407081
$data[0,4] = 16.86991607391245
$data[0,5] = 3.630603431827059
$data[0,6] = 0
$data[0,7] = 20.10074251947608
$data[0,8] = 0
$data[0,9] = 7.974778230901852
$data[0,10] = 0
$data[0,11] = 13.41094510034195
$data[0,12] = 18.31194102333822
$data[0,13] = 21.15023834395223
$data[1,0] = 7.786216191088958
$data[1,1] = 4.421841791944982
$data[1,2] = 0
$data[1,3] = 12.52275609309331
$data[1,4] = 15.89584955866815
$data[1,5] = 3.632132058032816
$data[1,6] = 0
$data[1,7] = 20.18974032786772
$data[1,8] = 0
$data[1,9] = 7.797541338779874
$data[1,10] = 0
$data[1,11] = 13.24434619055926
$data[1,12] = 18.36936698444557
$data[1,13] = 21.23208012869285
$data[2,0] = 7.612192862166771
$data[2,1] = 4.318149319971216
$data[2,2] = 0
$data[2,3] = 12.3892508566963
$data[2,4] = 15.26997757108489
$data[2,5] = 3.633120219987842
$data[2,6] = 0
$data[2,7] = 20.24804265683637
$data[2,8] = 0
$data[2,9] = 7.687811476672208
$data[2,10] = 0
$data[2,11] = 13.14417694576039
$data[2,12] = 18.40631605653857
$data[2,13] = 21.28670683551182
$data[3,0] = 7.540386589786907
$data[3,1] = 4.274920983179022
$data[3,2] = 0
$data[3,3] = 12.33554092165054
$data[3,4] = 15.008197319934
$data[3,5] = 3.633535409427155
$data[3,6] = 0
$data[3,7] = 20.27272088643588
$data[3,8] = 0
$data[3,9] = 7.642930967971266
$data[3,10] = 0
$data[3,11] = 13.10393904122187
$data[3,12] = 18.4217991185508
$data[3,13] = 21.31006607332367
$data[4,0] = 7.528412891145766
$data[4,1] = 4.267685360765378
$data[4,2] = 0
$data[4,3] = 12.32666671298655
$data[4,4] = 14.96433081551589
$data[4,5] = 3.633605107740895
$data[4,6] = 0
$data[4,7] = 20.27687422432457
$data[4,8] = 0
$data[4,9] = 7.635470532342256
$data[4,10] = 0
$data[4,11] = 13.09729409201086
$data[4,12] = 18.42439583560536
$data[4,13] = 21.31401114405286
$data[5,0] = 7.611227908038814
$data[5,1] = 4.317570214357035
$data[5,2] = 0
$data[5,3] = 12.38852358776245
$data[5,4] = 15.26647399323133
$data[5,5] = 3.633125768693229
$data[5,6] = 0
$data[5,7] = 20.2483717529377
$data[5,8] = 0
$data[5,9] = 7.687206784353235
$data[5,10] = 0
$data[5,11] = 13.14363186474681
$data[5,12] = 18.4065231399831
$data[5,13] = 21.28701742139987
$data[6,0] = 7.968554709572069
$data[6,1] = 4.529060539686277
$data[6,2] = 0
$data[6,3] = 12.66736579649584
$data[6,4] = 16.5399640634477
$data[6,5] = 3.631120235057661
$data[6,6] = 0
$data[6,7] = 20.1306698210983
$data[6,8] = 0
$data[6,9] = 7.913892027303334
$data[6,10] = 0
$data[6,11] = 13.3530902439452
$data[6,12] = 18.33139162178571
$data[6,13] = 21.17754814653876
$data[7,0] = 8.633265719741289
$data[7,1] = 4.910515248819764
$data[7,2] = 0
$data[7,3] = 13.22997969672598
$data[7,4] = 19.00274580682531
$data[7,5] = 3.627579027650773
$data[7,6] = 0
$data[7,7] = 19.92888051642592
$data[7,8] = 0
$data[7,9] = 8.348590763309826
$data[7,10] = 0
$data[7,11] = 13.77853092630695
$data[7,12] = 18.19740601033532
$data[7,13] = 20.99767911913446
$data[8,0] = 9.094483978192073
$data[8,1] = 5.16886182939241
$data[8,2] = 0
$data[8,3] = 13.6482231143898
$data[8,4] = 20.67494806633232
$data[8,5] = 3.625213599359153
$data[8,6] = 0
$data[8,7] = 19.79832611648771
$data[8,8] = 0
$data[8,9] = 8.658632704718439
$data[8,10] = 0
$data[8,11] = 14.09714194225444
$data[8,12] = 18.10702294116326
$data[8,13] = 20.88685574293936
$data[9,0] = 9.297498834443923
$data[9,1] = 5.281359815348089
$data[9,2] = 0
$data[9,3] = 13.83857806799188
$data[9,4] = 21.3917225636224
$data[9,5] = 3.624188291881028
$data[9,6] = 0
$data[9,7] = 19.74277922163877
$data[9,8] = 0
$data[9,9] = 8.797010234693699
$data[9,10] = 0
$data[9,11] = 14.24277833614572
$data[9,12] = 18.06763751087496
$data[9,13] = 20.84109380699267
$data[10,0] = 9.373333722953268
$data[10,1] = 5.323218379137062
$data[10,2] = 0
$data[10,3] = 13.91059590378359
$data[10,4] = 21.65686569030329
$data[10,5] = 3.623807291191609
$data[10,6] = 0
$data[10,7] = 19.72229805187377
$data[10,8] = 0
$data[10,9] = 8.848978513719345
$data[10,10] = 0
$data[10,11] = 14.29797644192935
$data[10,12] = 18.05297082128205
$data[10,13] = 20.82443590647283
$data[11,0] = 9.357048640362226
$data[11,1] = 5.314236668369573
$data[11,2] = 0
$data[11,3] = 13.89508970885303
$data[11,4] = 21.60004134736742
$data[11,5] = 3.623889024126606
$data[11,6] = 0
$data[11,7] = 19.72668442265932
$data[11,8] = 0
$data[11,9] = 8.837806191094517
$data[11,10] = 0
$data[11,11] = 14.28608718581526
$data[11,12] = 18.05611855527268
$data[11,13] = 20.82799359942724
$data[12,0] = 9.303759053867164
$data[12,1] = 5.284818516437577
$data[12,2] = 0
$data[12,3] = 13.84450477152222
$data[12,4] = 21.4136618050453
$data[12,5] = 3.624156801403152
$data[12,6] = 0
$data[12,7] = 19.74108313213342
$data[12,8] = 0
$data[12,9] = 8.801294630739671
$data[12,10] = 0
$data[12,11] = 14.24731883277698
$data[12,12] = 18.0664259150808
$data[12,13] = 20.83970988826042
$data[13,0] = 9.270980087058886
$data[13,1] = 5.266701892012538
$data[13,2] = 0
$data[13,3] = 13.81350922296411
$data[13,4] = 21.29868154950795
$data[13,5] = 3.624321767279518
$data[13,6] = 0
$data[13,7] = 19.74997482399189
$data[13,8] = 0
$data[13,9] = 8.778872541371619
$data[13,10] = 0
$data[13,11] = 14.2235769302419
$data[13,12] = 18.07277169851923
$data[13,13] = 20.84697392126693
$data[14,0] = 9.081073976313725
$data[14,1] = 5.161406984369471
$data[14,2] = 0
$data[14,3] = 13.63577804040474
$data[14,4] = 20.62722412089977
$data[14,5] = 3.625281623586128
$data[14,6] = 0
$data[14,7] = 19.80203361893024
$data[14,8] = 0
$data[14,9] = 8.649531528272565
$data[14,10] = 0
$data[14,11] = 14.08763387867796
$data[14,12] = 18.10963159552224
$data[14,13] = 20.8899401932099
$data[15,0] = 8.962783238981617
$data[15,1] = 5.095510248661035
$data[15,2] = 0
$data[15,3] = 13.52671263079648
$data[15,4] = 20.20408069597325
$data[15,5] = 3.625883434476691
$data[15,6] = 0
$data[15,7] = 19.83495470511188
$data[15,8] = 0
$data[15,9] = 8.569466842092913
$data[15,10] = 0
$data[15,11] = 14.00437908280663
$data[15,12] = 18.13268633254065
$data[15,13] = 20.91749169837002
$data[16,0] = 8.894109267074239
$data[16,1] = 5.057136171225134
$data[16,2] = 0
$data[16,3] = 13.46399447655786
$data[16,4] = 19.95656407809801
$data[16,5] = 3.626234357932931
$data[16,6] = 0
$data[16,7] = 19.85425172048357
$data[16,8] = 0
$data[16,9] = 8.523167916054323
$data[16,10] = 0
$data[16,11] = 13.95656192891812
$data[16,12] = 18.14610973437391
$data[16,13] = 20.93377635155587
$data[17,0] = 8.870750279150991
$data[17,1] = 5.044062933396289
$data[17,2] = 0
$data[17,3] = 13.44276394199352
$data[17,4] = 19.87204792380568
$data[17,5] = 3.626353996240438
$data[17,6] = 0
$data[17,7] = 19.86084745854642
$data[17,8] = 0
$data[17,9] = 8.507450895139421
$data[17,10] = 0
$data[17,11] = 13.94038526978751
$data[17,12] = 18.15068268398909
$data[17,13] = 20.93936517367581
$data[18,0] = 8.975441838990879
$data[18,1] = 5.102574069642547
$data[18,2] = 0
$data[18,3] = 13.53832202640121
$data[18,4] = 20.24955283636154
$data[18,5] = 3.625818876455775
$data[18,6] = 0
$data[18,7] = 19.83141276266693
$data[18,8] = 0
$data[18,9] = 8.578015884457697
$data[18,10] = 0
$data[18,11] = 14.01323495469956
$data[18,12] = 18.13021526206347
$data[18,13] = 20.91451347164319
$data[19,0] = 9.31944026925769
$data[19,1] = 5.293479619035115
$data[19,2] = 0
$data[19,3] = 13.85936518382207
$data[19,4] = 21.46857628470577
$data[19,5] = 3.624077951942392
$data[19,6] = 0
$data[19,7] = 19.73683886235736
$data[19,8] = 0
$data[19,9] = 8.812031069916163
$data[19,10] = 0
$data[19,11] = 14.25870512141089
$data[19,12] = 18.06339167931378
$data[19,13] = 20.83625029817295
$data[20,0] = 9.53816220854236
$data[20,1] = 5.413915441140591
$data[20,2] = 0
$data[20,3] = 14.06876961769047
$data[20,4] = 22.22866616901552
$data[20,5] = 3.622982463754532
$data[20,6] = 0
$data[20,7] = 19.67825434044349
$data[20,8] = 0
$data[20,9] = 8.962433375804324
$data[20,10] = 0
$data[20,11] = 14.41939508712673
$data[20,12] = 18.02116199643291
$data[20,13] = 20.78901376293333
$data[21,0] = 9.422003970884047
$data[21,1] = 5.35003867304509
$data[21,2] = 0
$data[21,3] = 13.95706982989551
$data[21,4] = 21.82633154458858
$data[21,5] = 3.62356328685494
$data[21,6] = 0
$data[21,7] = 19.70922669908903
$data[21,8] = 0
$data[21,9] = 8.882408684313935
$data[21,10] = 0
$data[21,11] = 14.33362475182234
$data[21,12] = 18.04356906094873
$data[21,13] = 20.81386604244074
$data[22,0] = 8.969720958601689
$data[22,1] = 5.099382038354141
$data[22,2] = 0
$data[22,3] = 13.53307345931894
$data[22,4] = 20.22900810905287
$data[22,5] = 3.625848047763973
$data[22,6] = 0
$data[22,7] = 19.83301292159118
$data[22,8] = 0
$data[22,9] = 8.574151695448059
$data[22,10] = 0
$data[22,11] = 14.00923106422871
$data[22,12] = 18.13133190673453
$data[22,13] = 20.91585854217891
$data[23,0] = 8.457871838811615
$data[23,1] = 4.811065505900966
$data[23,2] = 0
$data[23,3] = 13.07658009095259
$data[23,4] = 18.34778573295695
$data[23,5] = 3.628495345511567
$data[23,6] = 0
$data[23,7] = 19.98036216121421
$data[23,8] = 0
$data[23,9] = 8.232402515009033
$data[23,10] = 0
$data[23,11] = 13.66216564262268
$data[23,12] = 18.23223214622903
$data[23,13] = 21.04260131724284
$ws.Range("B2:O25").Value = $data
